# Add 2022-Q3 data:
#  - insert a new worksheet "2022-Q3" right after "总计" (duplicated from the
#    "2022-Q2" sheet, which has the same column layout), with refreshed figures
#  - "总计" (summary) sheet gets a new leading row for the 2022-Q3 entry, and
#    the previously-existing rows shift down by one

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as TEXT (inlineStr/shared-string), even when the text
# looks like a number (e.g. "0.77"), without leaving a stray quote-prefixed
# number-format style behind on the cell.
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range("ZZ9999").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 1) Duplicate the "2022-Q2" sheet to create the new "2022-Q3" sheet, placed
#    immediately before it (so the tab order becomes 总计, 2022-Q3, 2022-Q2,
#    2021-Q3), then refresh the fund's quarterly figures.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)

$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

Set-TextValue $wsQ3 "D2" "0.77"
Set-TextValue $wsQ3 "E2" "89.31"
Set-TextValue $wsQ3 "F2" "4.06"
Set-TextValue $wsQ3 "G2" "0.0313"
$wsQ3.Range("H2").Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: shift the existing two rows down by one
#    and insert the new 2022-Q3 row at the top of the data.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Row 4 becomes what row 3 used to hold ("2021-Q3"); copy A3's style for A4.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q3"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0

# Row 3 becomes what row 2 used to hold ("2022-Q2").
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.03

# Row 2 becomes the new "2022-Q3" entry (counts unchanged from old row 2).
$wsTotal.Range("B2").Value = "2022-Q3"
